$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B$($r1):AD$($r1)")
    $rng2 = $ws.Range("B$($r2):AD$($r2)")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# Simple two-row swaps (rank/A column stays put; id..PL_AhUnder columns B:AD swap)
Swap-Rows 95 96
Swap-Rows 110 111
Swap-Rows 129 131
Swap-Rows 200 201
Swap-Rows 224 225
Swap-Rows 231 232
Swap-Rows 237 238
Swap-Rows 249 250
Swap-Rows 256 257

# Three-way rotation: 267 <- 270 <- 268 <- 267 (old)
$arr267 = $ws.Range("B267:AD267").Value2
$arr268 = $ws.Range("B268:AD268").Value2
$arr270 = $ws.Range("B270:AD270").Value2
$ws.Range("B267:AD267").Value2 = $arr270
$ws.Range("B268:AD268").Value2 = $arr267
$ws.Range("B270:AD270").Value2 = $arr268
